# Fixing the "students" sheet and the "student_pswd" sheet: the last two
# demo rows (Kumara / `2h) were bogus test rows created while chasing the
# read_password bug - removing them and re-typing clean placeholder data
# for the remaining rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "students" sheet
# ---------------------------------------------------------------------
$students = $wb.Worksheets.Item("students")
$students.Activate()

$students.Range("C2").Value = "fdfdf"
$students.Range("D2").Value = "j"
$students.Range("E2").Value = "kj"
$students.Range("F2").Value = "klj"
$students.Range("G2").Value = "kl"
$students.Range("H2").Value = "jjkl"

$students.Range("B3").Value = "Lee"
$students.Range("C3").Value = "jkljk"
$students.Range("D3").Value = "jkl"
$students.Range("E3").Value = "jkl"
$students.Range("F3").Value = "jk"
$students.Range("G3").Value = "lj"
$students.Range("H3").Value = "kljkl"
$students.Range("J3").Value = 6

$students.Range("B4").Value = "Dave"
$students.Range("C4").Value = "hj"
$students.Range("D4").Value = "h"
$students.Range("E4").Value = "jkh"
$students.Range("H4").Value = "hjk"

$students.Range("B5").Value = "Dave"
$students.Range("C5").Value = "hj"
$students.Range("D5").Value = "hj"
$students.Range("E5").Value = "khjk"
$students.Range("F5").Value = "h"
$students.Range("G5").Value = "jkh"
$students.Range("H5").Value = "jkhjk"

$students.Range("B6").Value = "Sam"
$students.Range("C6").Value = "hjh"
$students.Range("D6").Value = "jkh"
$students.Range("E6").Value = "jkh"
$students.Range("F6").Value = "jk"
$students.Range("G6").Value = "hj"
$students.Range("H6").Value = "hjk"

$students.Range("B7").Value = "Agnetha"
$students.Range("C7").Value = "j"
$students.Range("D7").Value = "k"
$students.Range("E7").Value = "kl"
$students.Range("F7").Value = "kl"
$students.Range("G7").Value = "jk"
$students.Range("H7").Value = "lj"

# Remove the two trailing demo rows (Kumara, `2h) entirely.
[void]$students.Range("A8:A9").EntireRow.Delete()

[void]$students.Range("J3").Select()

# ---------------------------------------------------------------------
# "student_pswd" sheet
# ---------------------------------------------------------------------
$pswd = $wb.Worksheets.Item("student_pswd")
$pswd.Activate()

$pswd.Range("C3").Value = "Lee"
$pswd.Range("C4").Value = "Dave"
$pswd.Range("C5").Value = "Dave"
$pswd.Range("C6").Value = "Sam"
$pswd.Range("G6").Value = 6
$pswd.Range("C7").Value = "Agnetha"

# The old rows 8/9 (student7/Kumara, student8/`2h) are cleared out, not
# fully deleted - row 10 was already a blank formatted row beneath them.
[void]$pswd.Range("A8:C9").ClearContents()

[void]$pswd.Range("B2").Select()
